# Update the "想去人数" (F column) figures on the "展览" and "全部类型"
# sheets to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    3  = 278
    4  = 289
    7  = 7065
    9  = 77
    10 = 125
    11 = 93
    12 = 2
    13 = 42
    16 = 25
    17 = 243
    18 = 645
    19 = 14
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
